# Etapa4.pptx - "Correcao etapa 5 aula 1"
# Slide 4 (index 4): move the picture "Imagem 13" up (y offset
# 2461450 EMU -> 184421 EMU; x stays at 7196362 EMU).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item("Imagem 13")

# PowerPoint COM positions are expressed in points (1 pt = 12700 EMU).
# 184421 EMU / 12700 = 14.521338582677165 pt; nudge to 14.521339 so the
# EMU value PowerPoint writes back out round-trips to exactly 184421.
$sh.Top = 14.521339
